$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing three data rows (old rows 20, 21 and 22), which are
# no longer part of the accelerometer sample after the edit.
$ws.Rows("20:22").Delete()

# Insert two brand-new sample rows right after the header row; everything
# that used to start at row 2 shifts down to row 4.
$ws.Rows("2:3").Insert()

# Populate the two newly inserted rows with the new accelerometer readings.
$ws.Range("A2").Value = -3.747647881507874
$ws.Range("B2").Value = 4.277600646018982
$ws.Range("C2").Value = 0.2108629420399666

$ws.Range("A3").Value = -3.872398495674133
$ws.Range("B3").Value = 4.347799897193909
$ws.Range("C3").Value = 0.4252039864659314
